# Update countries & provincias Spain
# Applies numeric updates for several countries and reorders/updates three
# pairs of rows (Uruguay/Mali, Nueva Caledonia/Belice,
# Islas Virgenes Britanicas/Butan) to match the new source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4,2).Value = 1356629
$ws.Cells.Item(4,3).Value = 9320
$ws.Cells.Item(4,4).Value = 240616
$ws.Cells.Item(4,5).Value = 1035591
$ws.Cells.Item(4,6).Value = 16494
$ws.Cells.Item(4,7).Value = 385
$ws.Cells.Item(4,8).Value = 80422

# --- Row 10: Alemania ---
$ws.Cells.Item(10,2).Value = 171767
$ws.Cells.Item(10,3).Value = 443
$ws.Cells.Item(10,5).Value = 19810
$ws.Cells.Item(10,7).Value = 8
$ws.Cells.Item(10,8).Value = 7557

# --- Row 47: Sudafrica ---
$ws.Cells.Item(47,2).Value = 10015
$ws.Cells.Item(47,3).Value = 595
$ws.Cells.Item(47,4).Value = 4173
$ws.Cells.Item(47,5).Value = 5648
$ws.Cells.Item(47,7).Value = 8
$ws.Cells.Item(47,8).Value = 194

# --- Row 110: Burkina Faso ---
$ws.Cells.Item(110,2).Value = 751
$ws.Cells.Item(110,3).Value = 3
$ws.Cells.Item(110,4).Value = 577
$ws.Cells.Item(110,5).Value = 125
$ws.Cells.Item(110,7).Value = 1
$ws.Cells.Item(110,8).Value = 49

# --- Rows 113/114: Uruguay & Mali swap places, Mali gets new figures ---
$ws.Cells.Item(113,1).Value = "Mali"
$ws.Cells.Item(113,2).Value = 704
$ws.Cells.Item(113,3).Value = 12
$ws.Cells.Item(113,4).Value = 351
$ws.Cells.Item(113,5).Value = 315
$ws.Cells.Item(113,6).Value = 0
$ws.Cells.Item(113,7).Value = 1
$ws.Cells.Item(113,8).Value = 38

$ws.Cells.Item(114,1).Value = "Uruguay"
$ws.Cells.Item(114,2).Value = 702
$ws.Cells.Item(114,3).Value = 0
$ws.Cells.Item(114,4).Value = 513
$ws.Cells.Item(114,5).Value = 171
$ws.Cells.Item(114,6).Value = 8
$ws.Cells.Item(114,7).Value = 0
$ws.Cells.Item(114,8).Value = 18

# --- Rows 192/193: Nueva Caledonia & Belice swap places (data travels with name) ---
$ws.Cells.Item(192,1).Value = "Belice"
$ws.Cells.Item(192,4).Value = 16
$ws.Cells.Item(192,8).Value = 2

$ws.Cells.Item(193,1).Value = "Nueva Caledonia"
$ws.Cells.Item(193,4).Value = 18
$ws.Cells.Item(193,8).Value = 0

# --- Rows 212/213: Islas Virgenes Britanicas & Butan swap places ---
$ws.Cells.Item(212,1).Value = "Butan"
$ws.Cells.Item(212,4).Value = 5
$ws.Cells.Item(212,8).Value = 0

$ws.Cells.Item(213,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213,4).Value = 4
$ws.Cells.Item(213,8).Value = 1
